# "Generate Report for Handoff"
#
# Updates the localization-status report with a freshly generated
# handoff report:
#   - Bumps the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#     timestamps to reflect the new handoff generation run.
#   - Sets the "Priority" column to "ht" for the rows whose handoff has
#     just been (re)generated (rows 8, 9, 10, 11, 13, 14 - row 12 is not
#     part of this handoff run).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the "Latest Handoff Datetime" / generate-date timestamps.
$wsOverview.Cells.Replace("2016-08-28 20:22:18", "2016-08-28 20:22:45")
$wsDeDe.Cells.Replace("2016-08-28 20:22:18", "2016-08-28 20:22:45")
$wsZhCn.Cells.Replace("2016-08-28 20:22:13", "2016-08-28 20:22:39")

# Mark the freshly handed-off rows with priority "ht".
$rows = 8, 9, 10, 11, 13, 14
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}
